$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# NumberFormat "@" forces text storage so numeric-looking strings
# (e.g. "40.49", "1.50", "67.274.01") are preserved verbatim instead
# of being coerced into floating point numbers by COM; ClearFormats()
# afterwards restores the default (unstyled) cell formatting so only
# the cell content changes, matching the original workbook styling.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.274.01"
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +4.45%  "
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.257.53"
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +2.33%  "
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "577.84"
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +1.95%  "
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "178.95"
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +4.86%  "
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.602"
$cell.ClearFormats()

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -1.60%  "
$cell.ClearFormats()

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "3.256.34"
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +2.33%  "
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +3.64%  "
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.75"
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +1.56%  "
$cell.ClearFormats()

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +4.16%  "
$cell.ClearFormats()

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.830.19"
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +2.57%  "
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +0.79%  "
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "28.24"
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +2.79%  "
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "67.226.42"
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +4.36%  "
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +2.59%  "
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.259.35"
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  +2.42%  "
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.87"
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +1.85%  "
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.39"
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +2.76%  "
$cell.ClearFormats()

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "374.17"
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +5.72%  "
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "7.64"
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +5.88%  "
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "71.31"
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +3.10%  "
$cell.ClearFormats()

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.512"
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +1.42%  "
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.397.87"
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +2.60%  "
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.68%  "
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.86"
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +2.98%  "
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +1.44%  "
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +3.70%  "
$cell.ClearFormats()

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.62"
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "22.60"
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +2.18%  "
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +5.06%  "
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.81"
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +2.16%  "
$cell.ClearFormats()

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "164.31"
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +5.95%  "
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.50"
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  +3.88%  "
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.856"
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +4.19%  "
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.87"
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +9.43%  "
$cell.ClearFormats()

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "27.02"
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +4.10%  "
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +1.40%  "
$cell.ClearFormats()

$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value = "RenderToken"
$cell.ClearFormats()

$cell = $ws.Range("C43")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "6.58"
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +9.10%  "
$cell.ClearFormats()

$cell = $ws.Range("B44")
$cell.NumberFormat = "@"
$cell.Value = "Maker"
$cell.ClearFormats()

$cell = $ws.Range("C44")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.764.76"
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +5.27%  "
$cell.ClearFormats()

$cell = $ws.Range("B45")
$cell.NumberFormat = "@"
$cell.Value = "Filecoin"
$cell.ClearFormats()

$cell = $ws.Range("C45")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "4.40"
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +5.02%  "
$cell.ClearFormats()

$cell = $ws.Range("B46")
$cell.NumberFormat = "@"
$cell.Value = "InjectiveProtocol"
$cell.ClearFormats()

$cell = $ws.Range("C46")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell.ClearFormats()

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "25.60"
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +7.30%  "
$cell.ClearFormats()

$cell = $ws.Range("B47")
$cell.NumberFormat = "@"
$cell.Value = "Bittensor"
$cell.ClearFormats()

$cell = $ws.Range("C47")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "344.89"
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +6.09%  "
$cell.ClearFormats()

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "40.49"
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  +2.13%  "
$cell.ClearFormats()

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0674"
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +2.36%  "
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +3.21%  "
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +0.60%  "
$cell.ClearFormats()

